# Weekly update: a new daily price observation is inserted as row 30
# (pushing the existing rows 30-127 down to 31-128, i.e. the data keeps
# its relative order but everything from the old row 30 onward shifts
# down by one row).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before the current row 30; this shifts rows
# 30:127 down to 31:128 (formats/styles travel with the shift).
$ws.Rows("30:30").Insert()

# Populate the newly inserted row 30 with the new price observation.
$ws.Range("A30").Value = 1
$ws.Range("B30").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C30").Value = "Arica y Parinacota"
$ws.Range("D30").Value = 44953
$ws.Range("E30").Value = 15
$ws.Range("F30").Value = "Fruta"
$ws.Range("G30").Value = 100109
$ws.Range("H30").Value = "Uva"
$ws.Range("I30").Value = 100109001
$ws.Range("J30").Value = "Uva"
$ws.Range("K30").Value = "Red Globe"
$ws.Range("L30").Value = "Primera"
$ws.Range("M30").Value = 300
$ws.Range("N30").Value = 24000
$ws.Range("O30").Value = 25000
$ws.Range("P30").Value = 24500
$ws.Range("Q30").Value = "`$/caja 20 kilos"
$ws.Range("R30").Value = "Provincia de Limarí"
$ws.Range("S30").Value = 1225
$ws.Range("T30").Value = 20
